$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newIps = @("136.235.245.14","20.81.87.38","230.197.206.226","22.200.13.168","156.139.47.169","228.112.147.239","47.184.112.238","107.146.99.211","157.123.85.146","83.162.102.110","72.44.37.158","177.140.182.155","136.235.245.14","20.81.87.38","230.197.206.226","22.200.13.168","156.139.47.169","228.112.147.239","47.184.112.238","107.146.99.211","157.123.85.146","83.162.102.110","72.44.37.158","177.140.182.155","136.235.245.14","20.81.87.38","230.197.206.226","22.200.13.168","156.139.47.169","228.112.147.239","47.184.112.238","107.146.99.211","157.123.85.146","83.162.102.110","72.44.37.158","177.140.182.155","136.235.245.14","20.81.87.38","230.197.206.226","22.200.13.168","156.139.47.169","228.112.147.239","47.184.112.238","107.146.99.211","157.123.85.146","83.162.102.110","72.44.37.158","177.140.182.155","136.235.245.14","20.81.87.38","230.197.206.226","136.235.245.14","20.81.87.38","230.197.206.226","136.235.245.14","20.81.87.38","230.197.206.226","136.235.245.14","20.81.87.38","230.197.206.226","136.235.245.14","20.81.87.38","230.197.206.226","136.235.245.14","20.81.87.38","230.197.206.226")

$startRow = 27
for ($i = 0; $i -lt $newIps.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newIps[$i]
}
